$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.8168139999999999
$ws.Cells.Item(2, 8).Value = 2.450442
$ws.Cells.Item(2, 9).Value = 0.009738363985633989
$ws.Cells.Item(2, 10).Value = 0.00977835907772915
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1375686666666667
$ws.Cells.Item(2, 14).Value = 0.412706
$ws.Cells.Item(2, 15).Value = 0.2896572731203081
$ws.Cells.Item(2, 16).Value = 0.2896572731203081
$ws.Cells.Item(2, 17).Value = 0.1123680128946667
$ws.Cells.Item(2, 18).Value = 1.011312116052
$ws.Cells.Item(2, 19).Value = 0.002820787956731757
$ws.Cells.Item(2, 20).Value = 0.002832372826046237
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.8168139999999999
$ws.Cells.Item(3, 8).Value = 2.450442
$ws.Cells.Item(3, 9).Value = 0.009738363985633989
$ws.Cells.Item(3, 10).Value = 0.00977835907772915
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.3373673333333334
$ws.Cells.Item(3, 14).Value = 1.012102
$ws.Cells.Item(3, 15).Value = 0.7103427268796919
$ws.Cells.Item(3, 16).Value = 0.7103427268796919
$ws.Cells.Item(3, 17).Value = 0.2755663610093333
$ws.Cells.Item(3, 18).Value = 2.480097249084
$ws.Cells.Item(3, 19).Value = 0.006917576028902233
$ws.Cells.Item(3, 20).Value = 0.006945986251682914
$ws.Cells.Item(4, 9).Value = 0.9340146796604594
$ws.Cells.Item(4, 10).Value = 0.9378506425784978
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.1375686666666667
$ws.Cells.Item(4, 14).Value = 0.412706
$ws.Cells.Item(4, 15).Value = 0.2896572731203081
$ws.Cells.Item(4, 16).Value = 0.2896572731203081
$ws.Cells.Item(4, 17).Value = 10.77731061631311
$ws.Cells.Item(4, 18).Value = 96.995795546818
$ws.Cells.Item(4, 19).Value = 0.2705441451647868
$ws.Cells.Item(4, 20).Value = 0.2716552597234164
$ws.Cells.Item(5, 9).Value = 0.9340146796604594
$ws.Cells.Item(5, 10).Value = 0.9378506425784978
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.3373673333333334
$ws.Cells.Item(5, 14).Value = 1.012102
$ws.Cells.Item(5, 15).Value = 0.7103427268796919
$ws.Cells.Item(5, 16).Value = 0.7103427268796919
$ws.Cells.Item(5, 17).Value = 26.42980143102289
$ws.Cells.Item(5, 18).Value = 237.868212879206
$ws.Cells.Item(5, 19).Value = 0.6634705344956726
$ws.Cells.Item(5, 20).Value = 0.6661953828550814
$ws.Cells.Item(6, 7).Value = 2.332585666666667
$ws.Cells.Item(6, 8).Value = 6.997757
$ws.Cells.Item(6, 9).Value = 0.02780996438561621
$ws.Cells.Item(6, 10).Value = 0.02792417885617889
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1375686666666667
$ws.Cells.Item(6, 14).Value = 0.412706
$ws.Cells.Item(6, 15).Value = 0.2896572731203081
$ws.Cells.Item(6, 16).Value = 0.2896572731203081
$ws.Cells.Item(6, 17).Value = 0.3208907000491111
$ws.Cells.Item(6, 18).Value = 2.888016300442
$ws.Cells.Item(6, 19).Value = 0.008055358449510476
$ws.Cells.Item(6, 20).Value = 0.008088441501604543
$ws.Cells.Item(7, 7).Value = 2.332585666666667
$ws.Cells.Item(7, 8).Value = 6.997757
$ws.Cells.Item(7, 9).Value = 0.02780996438561621
$ws.Cells.Item(7, 10).Value = 0.02792417885617889
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.3373673333333334
$ws.Cells.Item(7, 14).Value = 1.012102
$ws.Cells.Item(7, 15).Value = 0.7103427268796919
$ws.Cells.Item(7, 16).Value = 0.7103427268796919
$ws.Cells.Item(7, 17).Value = 0.786938206134889
$ws.Cells.Item(7, 18).Value = 7.082443855214001
$ws.Cells.Item(7, 19).Value = 0.01975460593610573
$ws.Cells.Item(7, 20).Value = 0.01983573735457435
$ws.Cells.Item(8, 7).Value = 1.0291985
$ws.Cells.Item(8, 8).Value = 2.058397
$ws.Cells.Item(8, 9).Value = 0.01227049194365979
$ws.Cells.Item(8, 10).Value = 0.008213924259590904
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.1375686666666667
$ws.Cells.Item(8, 14).Value = 0.412706
$ws.Cells.Item(8, 15).Value = 0.2896572731203081
$ws.Cells.Item(8, 16).Value = 0.2896572731203081
$ws.Cells.Item(8, 17).Value = 0.1415854653803334
$ws.Cells.Item(8, 18).Value = 0.8495127922820002
$ws.Cells.Item(8, 19).Value = 0.003554237236245204
$ws.Cells.Item(8, 20).Value = 0.002379222902649847
$ws.Cells.Item(9, 7).Value = 1.0291985
$ws.Cells.Item(9, 8).Value = 2.058397
$ws.Cells.Item(9, 9).Value = 0.01227049194365979
$ws.Cells.Item(9, 10).Value = 0.008213924259590904
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.3373673333333334
$ws.Cells.Item(9, 14).Value = 1.012102
$ws.Cells.Item(9, 15).Value = 0.7103427268796919
$ws.Cells.Item(9, 16).Value = 0.7103427268796919
$ws.Cells.Item(9, 17).Value = 0.3472179534156667
$ws.Cells.Item(9, 18).Value = 2.083307720494
$ws.Cells.Item(9, 19).Value = 0.008716254707414585
$ws.Cells.Item(9, 20).Value = 0.005834701356941057
$ws.Cells.Item(10, 7).Value = 1.355979666666667
$ws.Cells.Item(10, 8).Value = 4.067939
$ws.Cells.Item(10, 9).Value = 0.01616650002463063
$ws.Cells.Item(10, 10).Value = 0.01623289522800313
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1375686666666667
$ws.Cells.Item(10, 14).Value = 0.412706
$ws.Cells.Item(10, 15).Value = 0.2896572731203081
$ws.Cells.Item(10, 16).Value = 0.2896572731203081
$ws.Cells.Item(10, 17).Value = 0.1865403147704444
$ws.Cells.Item(10, 18).Value = 1.678862832934
$ws.Cells.Item(10, 19).Value = 0.004682744313033904
$ws.Cells.Item(10, 20).Value = 0.00470197616659105
$ws.Cells.Item(11, 7).Value = 1.355979666666667
$ws.Cells.Item(11, 8).Value = 4.067939
$ws.Cells.Item(11, 9).Value = 0.01616650002463063
$ws.Cells.Item(11, 10).Value = 0.01623289522800313
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 0.3373673333333334
$ws.Cells.Item(11, 14).Value = 1.012102
$ws.Cells.Item(11, 15).Value = 0.7103427268796919
$ws.Cells.Item(11, 16).Value = 0.7103427268796919
$ws.Cells.Item(11, 17).Value = 0.4574632441975556
$ws.Cells.Item(11, 18).Value = 4.117169197778
$ws.Cells.Item(11, 19).Value = 0.01148375571159673
$ws.Cells.Item(11, 20).Value = 0.01153091906141208
